$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Style Guide")

# Fill in row 14 with the new "div" tag row (previously blank).
$ws.Range("A14").Value = "div"
$ws.Range("B14:J14").Value = "-"

# Copy the style of A14 from the style used on the other tag cells in column A
# (e.g. A13), matching the header-style formatting that the diff shows (s=13).
$ws.Range("A13").Copy()
$ws.Range("A14").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update the active selection to K14, matching the saved workbook state.
$ws.Range("K14").Select()
